$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.374.82"
$ws.Range("E2").Value = "  +4.22%  "
$ws.Range("D3").Value = "'3.460.46"
$ws.Range("E3").Value = "  +2.73%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'571.11"
$ws.Range("E5").Value = "  +2.45%  "
$ws.Range("D6").Value = "'185.00"
$ws.Range("E6").Value = "  +5.93%  "
$ws.Range("E7").Value = "  +2.17%  "
$ws.Range("D8").Value = "'3.453.27"
$ws.Range("E8").Value = "  +2.86%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  +7.40%  "
$ws.Range("D11").Value = "'0.647"
$ws.Range("E11").Value = "  +2.72%  "
$ws.Range("D12").Value = "'55.65"
$ws.Range("E12").Value = "  +3.42%  "
$ws.Range("E13").Value = "  +2.90%  "
$ws.Range("D14").Value = "'9.40"
$ws.Range("E14").Value = "  +3.74%  "
$ws.Range("D15").Value = "'4.007.88"
$ws.Range("E15").Value = "  +2.46%  "
$ws.Range("D16").Value = "'18.58"
$ws.Range("E16").Value = "  +2.03%  "
$ws.Range("D17").Value = "'3.462.63"
$ws.Range("E17").Value = "  +2.64%  "
$ws.Range("D18").Value = "'67.293.10"
$ws.Range("E18").Value = "  +4.15%  "
$ws.Range("E19").Value = "  +0.87%  "
$ws.Range("D20").Value = "'12.05"
$ws.Range("E20").Value = "  +2.67%  "
$ws.Range("E21").Value = "  +2.50%  "
$ws.Range("D22").Value = "'488.08"
$ws.Range("E22").Value = "  +6.23%  "
$ws.Range("E23").Value = "  +2.34%  "
$ws.Range("D24").Value = "'15.12"
$ws.Range("E24").Value = "  +11.37%  "
$ws.Range("D26").Value = "'90.29"
$ws.Range("E26").Value = "  +4.72%  "
$ws.Range("E27").Value = "  +0.48%  "
$ws.Range("D28").Value = "'10.99"
$ws.Range("E28").Value = "  +2.08%  "
$ws.Range("E29").Value = "  +2.89%  "
$ws.Range("D30").Value = "'31.63"
$ws.Range("E30").Value = "  +3.55%  "
$ws.Range("D31").Value = "'6.99"
$ws.Range("E31").Value = "  +4.88%  "
$ws.Range("D32").Value = "'593.86"
$ws.Range("E32").Value = "  +3.96%  "
$ws.Range("D33").Value = "'11.64"
$ws.Range("D34").Value = "'63.46"
$ws.Range("E34").Value = "  +3.85%  "
$ws.Range("E35").Value = "  +2.28%  "
$ws.Range("D36").Value = "'0.148"
$ws.Range("E36").Value = "  +6.71%  "
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").Value = "'3.66"
$ws.Range("E38").Value = "  +0.66%  "
$ws.Range("D39").Value = "'0.0₃0789"
$ws.Range("E39").Value = "  +6.65%  "
$ws.Range("D40").Value = "'0.389"
$ws.Range("E40").Value = "  +5.63%  "
$ws.Range("D41").Value = "'36.61"
$ws.Range("E41").Value = "  +3.57%  "
$ws.Range("D42").Value = "'3.141.86"
$ws.Range("E42").Value = "  +2.22%  "
$ws.Range("E43").Value = "  +3.81%  "
$ws.Range("E44").Value = "  +7.29%  "
$ws.Range("E45").Value = "  +2.74%  "
$ws.Range("E46").Value = "  +21.93%  "
$ws.Range("D47").Value = "'3.27"
$ws.Range("E47").Value = "  +5.01%  "
$ws.Range("E48").Value = "  +1.21%  "
$ws.Range("D49").Value = "'8.78"
$ws.Range("E49").Value = "  +7.13%  "
$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "'142.05"
$ws.Range("E51").Value = "  +2.70%  "
